# Made excel reader a generic method.
# Adds a "logout" test case (rows 7-10) mirroring the existing "Login" case,
# plus the start of a "signin" test case (rows 11-16), on the Keyword sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window / view changes -------------------------------------------------
$wb.Windows.Item(1).Left = 240
$wb.Windows.Item(1).Top = 240
$wb.Windows.Item(1).Width = 25360
$wb.Windows.Item(1).Height = 14000

# --- Remove the stray empty A3 cell (it disappears entirely, not just blanks) ---
$ws.Range("A3").Clear()

# --- New column marker on the header/table's last existing row (F6) --------
$ws.Range("F6").Borders.Item(10).LineStyle = 1
$ws.Range("F6").Borders.Item(7).LineStyle = 1

# --- "logout" test case (rows 7-10), modeled on the "Login" block above ----
$ws.Range("A7").Value = "logout"
$ws.Range("B8").Value = "GOTOURL"
$ws.Range("E8").Value = "url"
$ws.Range("B9").Value = "SETTEXT"
$ws.Range("C9").Value = "username"
$ws.Range("D9").Value = "id"
$ws.Range("E9").Value = "Admin"
$ws.Range("B10").Value = "SETTEXT"
$ws.Range("C10").Value = "password"
$ws.Range("D10").Value = "id"
$ws.Range("E10").Value = "admin"

$rng1 = $ws.Range("B7:E10")
$rng1.Font.Color = 0

$boxCells = @("B7","C7","D7","E7","B8","C8","D8","E8","B9","C9","D9","E9","B10","C10","D10","E10")
foreach ($c in $boxCells) {
    $ws.Range($c).Borders.Item(10).LineStyle = 1
    $ws.Range($c).Borders.Item(9).LineStyle = 1
}
foreach ($c in @("B7","B8","B9","B10")) {
    $ws.Range($c).Borders.Item(7).LineStyle = 1
}
foreach ($c in @("B7","C7","D7","E7")) {
    $ws.Range($c).Borders.Item(8).LineStyle = 1
}

# --- "signin" test case (rows 11-16), drafted from the "logout" block ------
$ws.Range("A11").Value = "signin"
$ws.Range("B12").Value = "signin"
$ws.Range("E12").Value = "url"
$ws.Range("B13").Value = "signin"
$ws.Range("C13").Value = "username"
$ws.Range("D13").Value = "signin"
$ws.Range("E13").Value = "signin"
$ws.Range("B14").Value = "signin"
$ws.Range("C14").Value = "signin"
$ws.Range("D14").Value = "id"
$ws.Range("E14").Value = "signin"
$ws.Range("D15").Value = "signin"
$ws.Range("C16").Value = "signin"

$rng2 = $ws.Range("B11:E14")
$rng2.Font.Color = 0

$boxCells2 = @("B11","C11","D11","E11","B12","C12","D12","E12","B13","C13","D13","E13","B14","C14","D14","E14")
foreach ($c in $boxCells2) {
    $ws.Range($c).Borders.Item(10).LineStyle = 1
    $ws.Range($c).Borders.Item(9).LineStyle = 1
}
foreach ($c in @("B11","B12","B13","B14")) {
    $ws.Range($c).Borders.Item(7).LineStyle = 1
}
foreach ($c in @("B11","C11","D11","E11")) {
    $ws.Range($c).Borders.Item(8).LineStyle = 1
}

# --- Sheet view: zoomed in, selection parked on the last-typed cell --------
$ws.Application.ActiveWindow.Zoom = 175
$ws.Range("C16").Select()
